$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.687.83'
$ws.Range('E2').Value = '  +1.51%  '

$ws.Range('D3').Value = '1.636.89'
$ws.Range('E3').Value = '  +1.01%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '

$ws.Range('E6').Value = '  +3.87%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +2.76%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0624'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.68%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.55%  '

$ws.Range('D12').Value = '1.864.16'
$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('D13').Value = '1.640.36'
$ws.Range('E13').Value = '  +1.37%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.59%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.96%  '

$ws.Range('D16').Value = '26.692.48'
$ws.Range('E16').Value = '  +1.47%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.75%  '

$ws.Range('D18').Value = '0.0₃0747'
$ws.Range('E18').Value = '  +2.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.25%  '

$ws.Range('E20').Value = '  +0.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.07%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.29%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.93%  '

$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  +1.43%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.53%  '

$ws.Range('E30').Value = '  -2.56%  '

$ws.Range('E31').Value = '  -0.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.44%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.24%  '

$ws.Range('E34').Value = '  +1.53%  '

$ws.Range('E35').Value = '  -0.71%  '

$ws.Range('D36').Value = '1.215.10'
$ws.Range('E36').Value = '  +2.83%  '

$ws.Range('E37').Value = '  +5.96%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.813'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.86%  '

$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.507'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.29'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.34%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.09%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.795'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.30%  '

$ws.Range('D44').Value = '1.772.66'
$ws.Range('E44').Value = '  +0.92%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.37%  '

$ws.Range('E46').Value = '  +2.84%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.04'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.16%  '

$ws.Range('E48').Value = '  +0.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.17%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.411'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.37%  '

$ws.Range('E51').Value = '  +0.09%  '
